$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.266.55'
$ws.Range('E2').Value = '  +1.07%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.41'
$ws.Range('E3').Value = '  +1.31%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.46'
$ws.Range('E5').Value = '  +0.96%  '

$ws.Range('E6').Value = '  +0.35%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5136'
$ws.Range('E7').Value = '  +0.45%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3917'
$ws.Range('E8').Value = '  +2.96%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08361'
$ws.Range('E9').Value = '  +0.85%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.122'
$ws.Range('E10').Value = '  +1.61%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.61'
$ws.Range('E11').Value = '  +0.77%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.235'
$ws.Range('E12').Value = '  +0.78%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.73'
$ws.Range('E13').Value = '  +1.59%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.871.37'
$ws.Range('E14').Value = '  +0.71%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.295'
$ws.Range('E15').Value = '  +1.72%  '

$ws.Range('E16').Value = '  +0.26%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001106'
$ws.Range('E17').Value = '  +1.22%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '91.50'
$ws.Range('E18').Value = '  +1.40%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06664'
$ws.Range('E19').Value = '  +0.75%  '

$ws.Range('E20').Value = '  +0.63%  '

$ws.Range('E21').Value = '  +0.37%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.055'
$ws.Range('E22').Value = '  +0.92%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.301.67'
$ws.Range('E23').Value = '  +1.02%  '

$ws.Range('E24').Value = '  +1.54%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.279'
$ws.Range('E25').Value = '  +0.95%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.098.45'
$ws.Range('E26').Value = '  +1.20%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.525'
$ws.Range('E27').Value = '  -1.50%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '159.36'
$ws.Range('E28').Value = '  +1.46%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '20.67'
$ws.Range('E29').Value = '  +1.38%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.69'
$ws.Range('E30').Value = '  +0.97%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1067'
$ws.Range('E31').Value = '  +0.92%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.049'
$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.897'
$ws.Range('E33').Value = '  +5.74%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.608'
$ws.Range('E34').Value = '  +0.10%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.761'
$ws.Range('E35').Value = '  +1.36%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02468'
$ws.Range('E36').Value = '  +2.61%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06583'
$ws.Range('E37').Value = '  +1.16%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2196'
$ws.Range('E38').Value = '  +2.10%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.213'
$ws.Range('E39').Value = '  +0.62%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6529'
$ws.Range('E40').Value = '  +2.24%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.030'
$ws.Range('E41').Value = '  +3.42%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.233'
$ws.Range('E42').Value = '  +0.24%  '

$ws.Range('E43').Value = '  +0.87%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6159'
$ws.Range('E44').Value = '  +1.93%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.21'
$ws.Range('E45').Value = '  +0.94%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.289'
$ws.Range('E46').Value = '  +0.55%  '

$ws.Range('E47').Value = '  +0.85%  '

$ws.Range('E48').Value = '  +2.57%  '

$ws.Range('E49').Value = '  +2.73%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '121.59'
$ws.Range('E50').Value = '  +1.02%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.13'
$ws.Range('E51').Value = '  -0.53%  '
